$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the previously-missing "user story编号" values for rows 4 and 5
$ws.Range("D4").Value = 11
$ws.Range("D5").Value = 11

# New row 6: 封禁用户 (block users)
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 11.19
$ws.Range("C6").Value = "宋壬初 雷建坤"
$ws.Range("D6").Value = 20
$ws.Range("E6").Value = "封禁用户"

# New row 7: 邀请好友加入小组 (invite friends to join group)
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 11.19
$ws.Range("C7").Value = "郭俊石"
$ws.Range("D7").Value = 9
$ws.Range("E7").Value = "邀请好友加入小组"

$ws.Range("D8").Select()
